# input data update with Invalid Login data
#
# Adds a new "InvalidLogin" worksheet (after the existing "ValidLogin"
# sheet) containing a sample invalid-login row, mirroring the structure
# of the ValidLogin sheet, and makes it the active tab.

$wb = $excel.ActiveWorkbook
$validSheet = $wb.Worksheets.Item(1)

# Insert the new sheet right after ValidLogin so it becomes sheet index 2
# / the active tab, matching "activeTab=1" in the saved workbook.
$invalidSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $validSheet)
$invalidSheet.Name = "InvalidLogin"

# Header row
$invalidSheet.Range("A1").Value = "Username"
$invalidSheet.Range("B1").Value = "Password"

# Data row
$invalidSheet.Range("A2").Value = "meeralnissa123@gmail"
$invalidSheet.Range("B2").Value = "Actime"

# The username cell carries a mailto hyperlink (same pattern as ValidLogin!A2)
$invalidSheet.Hyperlinks.Add($invalidSheet.Range("A2"), "mailto:meeralnissa123@gmail", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "meeralnissa123@gmail")
$invalidSheet.Range("A2").ClearFormats()

$wb.Save()
